# Update the answer key table: each data row in the single table holds
# five "NNN÷N=Q, R" cells. This commit replaces the division problems
# with a new set of values (including one row where a cell was removed
# and a different cell appended, which nets out to the same 5-cell shape
# but with a full set of new cell contents).
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Map of (row, cell) -> new text, addressed positionally so there is no
# ambiguity from values that coincide between old/new cells in the same
# row (e.g. row 17 where "471÷8=58, 7" is both an old and a new value).
$updates = @(
    @{ Row = 1;  Cell = 1; Text = "395÷3=131, 2" },
    @{ Row = 1;  Cell = 2; Text = "749÷3=249, 2" },
    @{ Row = 1;  Cell = 3; Text = "336÷5=67, 1" },
    @{ Row = 1;  Cell = 4; Text = "685÷9=76, 1" },
    @{ Row = 1;  Cell = 5; Text = "405÷5=81, 0" },

    @{ Row = 5;  Cell = 1; Text = "583÷5=116, 3" },
    @{ Row = 5;  Cell = 2; Text = "444÷3=148, 0" },
    @{ Row = 5;  Cell = 3; Text = "747÷4=186, 3" },
    @{ Row = 5;  Cell = 4; Text = "695÷4=173, 3" },
    @{ Row = 5;  Cell = 5; Text = "254÷9=28, 2" },

    @{ Row = 9;  Cell = 1; Text = "585÷4=146, 1" },
    @{ Row = 9;  Cell = 2; Text = "484÷3=161, 1" },
    @{ Row = 9;  Cell = 3; Text = "292÷2=146, 0" },
    @{ Row = 9;  Cell = 4; Text = "675÷2=337, 1" },
    @{ Row = 9;  Cell = 5; Text = "125÷4=31, 1" },

    @{ Row = 13; Cell = 1; Text = "902÷4=225, 2" },
    @{ Row = 13; Cell = 2; Text = "877÷2=438, 1" },
    @{ Row = 13; Cell = 3; Text = "137÷9=15, 2" },
    @{ Row = 13; Cell = 4; Text = "685÷9=76, 1" },
    @{ Row = 13; Cell = 5; Text = "966÷4=241, 2" },

    @{ Row = 17; Cell = 1; Text = "292÷7=41, 5" },
    @{ Row = 17; Cell = 2; Text = "140÷9=15, 5" },
    @{ Row = 17; Cell = 3; Text = "471÷8=58, 7" },
    @{ Row = 17; Cell = 4; Text = "619÷5=123, 4" },
    @{ Row = 17; Cell = 5; Text = "235÷8=29, 3" }
)

foreach ($u in $updates) {
    $cell = $tbl.Rows.Item($u.Row).Cells.Item($u.Cell)
    $cell.Range.Text = $u.Text
}
